$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 337
    $ws.Range("F3").Value = 88
    $ws.Range("F4").Value = 1514
    $ws.Range("F6").Value = 44
}
